$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the sequential Day_Number / Date pattern into rows 43 and 44
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 43507

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 43508

# Copy the formatting from the previous data row so the new rows match
$ws.Range("A42:B42").Copy() | Out-Null
$ws.Range("A43:B44").PasteSpecial(-4122) | Out-Null

# Restore the values (PasteSpecial formats may have touched them) and reapply
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 43507
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 43508

# Keep the new rows' height as a non-custom (auto) height, matching the rest of the sheet
$ws.Range("A43:B44").EntireRow.AutoFit() | Out-Null

# Update the selection to match the new active cell/selection
$ws.Range("A43:A44").Select() | Out-Null
